# Update 1.6.3: Added new parameter for Project Progress
#
# - Remove two now-obsolete cell comments from "Материалы" sheet (G1, L1)
# - Replace the "Запланированный бюджет" header (L1) with a new
#   "Списано" column, reusing the yellow header style from the
#   neighbouring K1 cell
# - Switch sheet page orientation to portrait
# - Re-point the active sheet / selections back to "Материалы"

$wb = $excel.ActiveWorkbook
$wsMaterials = $wb.Worksheets.Item(1)
$wsServices  = $wb.Worksheets.Item(2)

# Drop the two stale comments on the "Материалы" sheet.
[void]$wsMaterials.Range("G1").Comment.Delete()
[void]$wsMaterials.Range("L1").Comment.Delete()

# Turn the old "Запланированный бюджет" column into "Списано",
# carrying over K1's header formatting (yellow fill) instead of the
# old blue one.
[void]$wsMaterials.Range("K1").Copy()
[void]$wsMaterials.Range("L1").PasteSpecial(-4122)
$wsMaterials.Range("L1").Value = "Списано"

# Switch the sheet to portrait orientation.
$wsMaterials.PageSetup.Orientation = 1

# Restore selections: "Услуги" back to D1 (no longer the active tab),
# then "Материалы" activated with I8 selected.
[void]$wsServices.Activate()
[void]$wsServices.Range("D1").Select()

[void]$wsMaterials.Activate()
[void]$wsMaterials.Range("I8").Select()
